# "first attempt in multiple linear regression model"
#
# The analyst wants a pristine copy of the genealogy matrix preserved on a
# second tab (to use as the full reference dataset) and then blanks out the
# last three animals (rows 12-14: 59K, 98K, 129K) on the original sheet,
# presumably to use them as a hold-out/test slice for the regression.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- 1. Duplicate Sheet1 -> Sheet2 --------------------------------------
# Worksheet.Copy(After:=ws1) gives us an exact snapshot of Sheet1 (values,
# shared-string references, number formats, row heights, column widths)
# positioned right after Sheet1, which is exactly the "before" state that
# belongs on the new tab.
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Sheet2"

# Sheet1 must stay the active/selected tab (Copy() makes the new sheet the
# active one by default).
$ws1.Activate()

# --- 2. Clear rows 12-14 on Sheet1 --------------------------------------
# Remove the data values but keep the cells (and their styles) in place.
$ws1.Range("A12:J14").ClearContents() | Out-Null

# The cleared rows go back to the sheet's default row height (15pt)
# instead of the data rows' custom 15.75pt height.
$ws1.Range("A12:J14").EntireRow.AutoFit() | Out-Null

# --- 3. Restore the selections described by the edit --------------------
$ws2.Range("B36").Select() | Out-Null
$ws1.Activate()
$ws1.Range("B14").Select() | Out-Null
